$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Folder Inventory": a freshly re-touched folder ("Azure Virtual
# Machine And Compute") moves back to the top of the list with a new
# "Last Updated" timestamp. Every row between the old and new position
# shifts down by one.
# ---------------------------------------------------------------------------
$wsInventory = $wb.Worksheets.Item("Folder Inventory")

$wsInventory.Range("A2").Value = "Azure Virtual Machine And Compute"
$wsInventory.Range("B2").Value = "Azure Virtual Machine And Compute"
$wsInventory.Range("C2").Value = "2025-06-12 17:31:45 +0530"

$zeroWidthSpace = [char]0x200B
$fabricPipelineTitle = "Work with Data Lake and Data Factory Pipelines in Microsoft Fabric" + $zeroWidthSpace
$wsInventory.Range("A3").Value = $fabricPipelineTitle
$wsInventory.Range("B3").Value = $fabricPipelineTitle
$wsInventory.Range("C3").Value = "2025-06-12 17:26:19 +0530"

$wsInventory.Range("A4").Value = "Get Started with Microsoft Fabric with Its Lakehouses"
$wsInventory.Range("B4").Value = "Get Started with Microsoft Fabric with Its Lakehouses"
$wsInventory.Range("C4").Value = "2025-06-12 16:16:30 +0530"

$wsInventory.Range("A5").Value = "Build A Fabric Real-Time Intelligence Solution in a Day"
$wsInventory.Range("B5").Value = "Build A Fabric Real-Time Intelligence Solution in a Day"
$wsInventory.Range("C5").Value = "2025-06-12 15:59:35 +0530"

# ---------------------------------------------------------------------------
# Sheet "Metadata": refresh generation timestamp and workflow run number.
# ---------------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B3").Value = "2025-06-12 12:04:26 UTC"

# "Workflow Run" holds a numeric-looking value but must stay a text cell
# (matching the source file, where it is stored as inline string "13").
# A leading apostrophe forces Excel to treat it as text; resetting the
# style afterwards drops the quote-prefix formatting Excel applies.
$wsMetadata.Range("B5").Value = "'13"
$wsMetadata.Range("B5").Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "Summary": most recent update timestamp matches the new top entry.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-12 17:31:45 +0530"
